# Update the "AddCustomerTest" sheet's runmode column (E) for the last two
# test rows (Ishita and Rohit) from "Y" to "N", and leave the sheet's
# selection on cell E6 (matching the final saved view state).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("AddCustomerTest")
$ws.Activate()

$ws.Range("E4").Value = "N"
$ws.Range("E5").Value = "N"

$ws.Range("E6").Select()
